$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: Digikey Order column (I) flips from "N" to "Y" ---
$ws.Range("I14").Value = "Y"

# --- Row 32: resistor part swapped to a new vendor/part (Vishay BC Components) ---
# Order matters for new shared-string allocation (K, D, H, J) so it matches the
# authored workbook's shared-strings ordering.
$ws.Range("K32").Value = "PPC2.49KXCT-ND"
$ws.Range("D32").Value = "RES 2.49K OHM 1/2W 1% AXIAL"
$ws.Range("H32").Value = "Vishay BC Components"
$ws.Range("J32").Value = "SFR16S0002491FR500"

# Row 32 shrinks from a 2-line wrapped description to a 1-line one.
$ws.Rows.Item(32).RowHeight = 27

# --- Row 33: unit price drop for the 3.9k 0.1% resistor ---
$ws.Range("L33").Value = 0.3

# --- Apply an AutoFilter on column I ("Digikey Order") showing only "N" rows ---
$ws.Range("A1:P47").AutoFilter(9, "N", 7)

# --- Selection moves off-sheet to D54 with no pinned top-left cell ---
$ws.Range("D54").Select()
